$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (also updates the tab name shown in the workbook)
$ws.Name = "Through 2022-11-07"

# Update the header label for the "through" date
$ws.Range("I1").Value = "2022 (through 11-07)"

# Update the data values for November (row 12) and Total (row 14)
$ws.Range("I12").Value = 18
$ws.Range("I14").Value = 1418
